# Auto-generated Excel COM-interop script.
#
# The workbook tracks Final Fantasy XIV "Leve" crafting/gathering profitability
# across 8 disciplines (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). A scheduled
# runner periodically refreshes the live market-board columns:
#   H = currentAveragePrice       I = currentAveragePriceNQ
#   J = currentAveragePriceHQ     K = LevePriceNQ
#   L = LevePriceHQ               M = LeveProfitNQ
#   N = LeveProfitHQ
# This script re-applies that refresh snapshot cell-by-cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
# Row 5
$ws.Range("H5").Value = 334.85715
$ws.Range("I5").Value = 57.333332
$ws.Range("K5").Value = 57.333332
$ws.Range("M5").Value = 57.666668
# Row 33
$ws.Range("H33").Value = 2820.3076
$ws.Range("I33").Value = 2779.3914
$ws.Range("J33").Value = 3134
$ws.Range("K33").Value = 2779.3914
$ws.Range("L33").Value = 3134
$ws.Range("M33").Value = -2550.3914
$ws.Range("N33").Value = -3592
# Row 86
$ws.Range("H86").Value = 2227071
$ws.Range("I86").Value = 4003928
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4003928
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -4002805
$ws.Range("N86").Value = -8246
# Row 89
$ws.Range("H89").Value = 2227071
$ws.Range("I89").Value = 4003928
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 20019640
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -20014024
$ws.Range("N89").Value = -41232
# Row 100
$ws.Range("H100").Value = 1840.375
$ws.Range("I100").Value = 942
$ws.Range("J100").Value = 2379.4
$ws.Range("K100").Value = 942
$ws.Range("L100").Value = 2379.4
$ws.Range("M100").Value = -401
$ws.Range("N100").Value = -3461.4
# Row 103
$ws.Range("H103").Value = 4953.1665
$ws.Range("J103").Value = 4953.1665
$ws.Range("L103").Value = 14859.4995
$ws.Range("N103").Value = -16031.4995
# Row 128
$ws.Range("H128").Value = 71450
$ws.Range("J128").Value = 71450
$ws.Range("L128").Value = 71450
$ws.Range("N128").Value = -81410
# Row 132
$ws.Range("H132").Value = 2796.125
$ws.Range("I132").Value = 2698
$ws.Range("K132").Value = 8094
$ws.Range("M132").Value = -5564
# Row 134
$ws.Range("H134").Value = 76666.336
$ws.Range("J134").Value = 76666.336
$ws.Range("L134").Value = 76666.336
$ws.Range("N134").Value = -86806.336
# Row 137
$ws.Range("H137").Value = 8903.6
$ws.Range("J137").Value = 11596.3
$ws.Range("L137").Value = 34788.89999999999
$ws.Range("N137").Value = -39888.89999999999
# Row 138
$ws.Range("H138").Value = 2673.795
$ws.Range("J138").Value = 3002.5693
$ws.Range("L138").Value = 9007.707900000001
$ws.Range("N138").Value = -19287.7079
# Row 141
$ws.Range("H141").Value = 4376.185
$ws.Range("I141").Value = 4294.32
$ws.Range("J141").Value = 5399.5
$ws.Range("K141").Value = 12882.96
$ws.Range("L141").Value = 16198.5
$ws.Range("M141").Value = -7702.959999999999
$ws.Range("N141").Value = -26558.5

# ---------------------------------------------------------------------------
# Sheet 2: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
# Row 46
$ws.Range("H46").Value = 56730.4
$ws.Range("I46").Value = 49499.668
$ws.Range("K46").Value = 49499.668
$ws.Range("M46").Value = -49180.668
# Row 61
$ws.Range("H61").Value = 12503992
$ws.Range("I61").Value = 9093984
$ws.Range("K61").Value = 9093984
$ws.Range("M61").Value = -9093772
# Row 74
$ws.Range("H74").Value = 11826296
$ws.Range("I74").Value = 15629350
$ws.Range("J74").Value = 1684817
$ws.Range("K74").Value = 15629350
$ws.Range("L74").Value = 1684817
$ws.Range("M74").Value = -15628476
$ws.Range("N74").Value = -1686565
# Row 77
$ws.Range("H77").Value = 11826296
$ws.Range("I77").Value = 15629350
$ws.Range("J77").Value = 1684817
$ws.Range("K77").Value = 78146750
$ws.Range("L77").Value = 8424085
$ws.Range("M77").Value = -78142382
$ws.Range("N77").Value = -8432821
# Row 136
$ws.Range("H136").Value = 12503992
$ws.Range("I136").Value = 9093984
$ws.Range("K136").Value = 27281952
$ws.Range("M136").Value = -27279402

# ---------------------------------------------------------------------------
# Sheet 3: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
# Row 20
$ws.Range("H20").Value = 5290.4
$ws.Range("I20").Value = 6067
$ws.Range("J20").Value = 4125.5
$ws.Range("K20").Value = 6067
$ws.Range("L20").Value = 4125.5
$ws.Range("M20").Value = -5820
$ws.Range("N20").Value = -4619.5
# Row 99
$ws.Range("H99").Value = 7066.5386
$ws.Range("I99").Value = 10217.077
$ws.Range("J99").Value = 3916
$ws.Range("K99").Value = 10217.077
$ws.Range("L99").Value = 3916
$ws.Range("M99").Value = -8719.076999999999
$ws.Range("N99").Value = -6912

# ---------------------------------------------------------------------------
# Sheet 4: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
# Row 19
$ws.Range("H19").Value = 7610.5
$ws.Range("I19").Value = 4026.25
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 4026.25
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -3856.25
$ws.Range("N19").Value = -10340
# Row 24
$ws.Range("H24").Value = 7610.5
$ws.Range("I24").Value = 4026.25
$ws.Range("J24").Value = 10000
$ws.Range("K24").Value = 4026.25
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = -3856.25
$ws.Range("N24").Value = -10340
# Row 31
$ws.Range("H31").Value = 273710.44
$ws.Range("I31").Value = 13631.272
$ws.Range("J31").Value = 424282.56
$ws.Range("K31").Value = 13631.272
$ws.Range("L31").Value = 424282.56
$ws.Range("M31").Value = -13336.272
$ws.Range("N31").Value = -424872.56
# Row 34
$ws.Range("H34").Value = 273710.44
$ws.Range("I34").Value = 13631.272
$ws.Range("J34").Value = 424282.56
$ws.Range("K34").Value = 13631.272
$ws.Range("L34").Value = 424282.56
$ws.Range("M34").Value = -13429.272
$ws.Range("N34").Value = -424686.56
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
# Row 132
$ws.Range("H132").Value = 1706.6364
$ws.Range("I132").Value = 1719.4062
$ws.Range("K132").Value = 5158.2186
$ws.Range("M132").Value = -2628.2186

# ---------------------------------------------------------------------------
# Sheet 5: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
# Row 2
$ws.Range("H2").Value = 93.42856999999999
$ws.Range("I2").Value = 66.36
$ws.Range("J2").Value = 115.258064
$ws.Range("K2").Value = 398.16
$ws.Range("L2").Value = 691.5483840000001
$ws.Range("M2").Value = -285.16
$ws.Range("N2").Value = -917.5483840000001
# Row 14
$ws.Range("H14").Value = 441.22223
$ws.Range("I14").Value = 441.22223
$ws.Range("K14").Value = 1323.66669
$ws.Range("M14").Value = -1150.66669
# Row 126
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 6000
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 18000
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -13060
$ws.Range("N126").Value = -39880
# Row 129
$ws.Range("H129").Value = 37042856
$ws.Range("I129").Value = 8672
$ws.Range("J129").Value = 83335580
$ws.Range("K129").Value = 26016
$ws.Range("L129").Value = 250006740
$ws.Range("M129").Value = -21016
$ws.Range("N129").Value = -250016740

# ---------------------------------------------------------------------------
# Sheet 6: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Range("H70").Value = 5974.75
$ws.Range("I70").Value = 5949.5
$ws.Range("K70").Value = 5949.5
$ws.Range("M70").Value = -5679.5
# Row 73
$ws.Range("H73").Value = 5974.75
$ws.Range("I73").Value = 5949.5
$ws.Range("K73").Value = 5949.5
$ws.Range("M73").Value = -5013.5
# Row 80
$ws.Range("H80").Value = 3006.8572
$ws.Range("I80").Value = 2976
$ws.Range("J80").Value = 3048
$ws.Range("K80").Value = 2976
$ws.Range("L80").Value = 3048
$ws.Range("M80").Value = -1978
$ws.Range("N80").Value = -5044
# Row 83
$ws.Range("H83").Value = 3006.8572
$ws.Range("I83").Value = 2976
$ws.Range("J83").Value = 3048
$ws.Range("K83").Value = 14880
$ws.Range("L83").Value = 15240
$ws.Range("M83").Value = -9888
$ws.Range("N83").Value = -25224

# ---------------------------------------------------------------------------
# Sheet 7: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(7)
# Row 55
$ws.Range("H55").Value = 34483180
$ws.Range("J55").Value = 488.5
$ws.Range("L55").Value = 488.5
$ws.Range("N55").Value = -834.5
# Row 132
$ws.Range("H132").Value = 159864.75
$ws.Range("I132").Value = 3751.4546
$ws.Range("J132").Value = 1113890.5
$ws.Range("K132").Value = 11254.3638
$ws.Range("L132").Value = 3341671.5
$ws.Range("M132").Value = -8724.363799999999
$ws.Range("N132").Value = -3346731.5
# Row 140
$ws.Range("H140").Value = 88819.25
$ws.Range("J140").Value = 88819.25
$ws.Range("L140").Value = 88819.25
$ws.Range("N140").Value = -99179.25

# ---------------------------------------------------------------------------
# Sheet 8: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(8)
# Row 122
$ws.Range("H122").Value = 1131
$ws.Range("I122").Value = 1123.7273
$ws.Range("K122").Value = 3371.1819
$ws.Range("M122").Value = -921.1819

